$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the format of the last-row (row 25, which has the bottom-border treatment)
# to row 20 (which will become the new final data row after MARTHA..JOHANA deletion).
$ws.Range("B25:J25").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)  # xlPasteFormats = -4122
$excel.CutCopyMode = 0

# Delete the 5 rows for MARTHA, KARELIS, CAROLINA, ANGELICA, JOHANA (old rows 21-25)
$ws.Rows("21:25").Delete()

# Update summary figures: Valor Mora total halved (10 -> 5 workers), Cant. Trabajadores 10 -> 5
$ws.Range("E11").Value = 284700
$ws.Range("C13").Value = 5

# Update Periodo Mora (column E) for all remaining worker rows to 2508
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"
$ws.Range("E19").Value = "2508"
$ws.Range("E20").Value = "2508"

# Row 18 previously held ELY CARMEN PRIMERA JULIO -> now becomes WENDY LILI TUÑON ARROYO
$ws.Range("C18").Value = "45554406"
$ws.Range("D18").Value = "WENDY LILI TUÑON ARROYO"

# Row 19 already held LORAINE VANESSA AVENDAÑO RIVERA's old data at row 20; reaffirm row 19 values (unchanged identity)
$ws.Range("C19").Value = "1143381535"
$ws.Range("D19").Value = "LORAINE VANESSA AVENDAÑO RIVERA"

# Row 20 previously held LORAINE's data -> now becomes the new worker MARIA DE LOS ANGELES LADEUX RODRIGUEZ
$ws.Range("C20").Value = "1007975454"
$ws.Range("D20").Value = "MARIA DE LOS ANGELES LADEUX RODRIGUEZ"

Write-Host "done"
